# SCADA/RemoteRuntime/tagsDeAplicacion.xlsx
# "agregado indicadores de valvulas y compuertas"
#
# Adds a new "Registered Tag List" row for a new internal tag,
# l21BtnLiberarAutManEnabled (enable/disable the LIBERAR AUTO/MAN button),
# inserted right after the existing l21BtnAutomaticoEnabled row, and flips
# a handful of existing tag default Values to reflect the new
# valve/gate-indicator wiring.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registered Tag List")

# --- Insert the new tag row (old row 40 "l21BtnInicioEnabled" and everything
#     below it shifts down by one, to row 41+) ---
$ws.Rows.Item(40).Insert()

$ws.Range("A40").Value = "Internal"
$ws.Range("B40").Value = "l21BtnLiberarAutManEnabled"
$ws.Range("C40").Value = "Boolean"
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = "Not Historical"
$ws.Range("H40").Value = "activa o desactiva el boton de LIBERAR AUTO/MAN"

# --- Update default Values on existing tags ---

# L21_TIEMPO_VACIADO / L21_TIEMPO_BYPASS and their textbox mirrors
$ws.Range("D27").Value = 10
$ws.Range("D28").Value = 10
$ws.Range("D30").Value = 10
$ws.Range("D31").Value = 10

# l21BtnPuestaEnMarchaEnabled / l21BtnResetEnabled / l21BtnManualEnabled /
# l21BtnAutomaticoEnabled now default enabled
$ws.Range("D36").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("D38").Value = 1
$ws.Range("D39").Value = 1

# l21CmbxOrigenValue / l21CmbxDestinoValue default selection (rows shifted
# down by the new row 40 insertion: old 49/51 -> new 50/52)
$ws.Range("D50").Value = 1
$ws.Range("D52").Value = 1

# l21TxtCantidadOrigenDestino default quantity (old row 54 -> new row 55)
$ws.Range("D55").Value = 13213
